$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (Excel would otherwise
# auto-convert them to a numeric type). We temporarily force text format,
# assign the value, then restore the original (default/Normal) style so the
# cell formatting matches the source workbook.
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D43", "D45", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply all cell value updates
$ws.Range("D2").Value = "26.102.67"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.655.02"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "217.82"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "0.5255"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.2610"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "0.06358"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").Value = "20.42"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").Value = "0.07796"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "4.513"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").Value = "1.620.28"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").Value = "0.5487"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "0.0₅8223"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "65.42"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "26.121.79"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D19").Value = "4.595"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "191.44"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "10.07"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "6.027"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "141.89"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("D26").Value = "7.253"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "1.428"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").Value = "0.05907"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "1.275"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").Value = "3.259"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "1.596"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "0.9521"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").Value = "2.784"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").Value = "2.409"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").Value = "0.5703"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").Value = "0.01621"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").Value = "5.813"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("D40").Value = "0.8492"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "1.030.76"
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").Value = "102.89"
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("D44").Value = "1.798.32"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "57.18"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").Value = "1.477"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.846"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05152"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "0.09702"
$ws.Range("E51").Value = "  +0.07%  "

# Restore default style on the cells we temporarily formatted as text
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
